$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Production table rows 87-91 ("produccion por linea de bebidas"):
#    the second-quarter quantity factor (columns H,K,N,Q) goes from 1 to 2,
#    and a couple of cells in column H gain explicit formulas.
# ---------------------------------------------------------------------------

# Row 87
$ws.Range("H87").Formula = "=E87+0"
$ws.Range("N87").Value = 2
$ws.Range("Q87").Value = 2

# Row 88 (master of a shared-formula block H88:H90)
$ws.Range("H88:H90").Formula = "=E88+0"
$ws.Range("Q88").Value = 2

# Row 89
$ws.Range("K89").Value = 2
$ws.Range("N89").Value = 2
$ws.Range("Q89").Value = 2

# Row 90
$ws.Range("E90").Value = 2
$ws.Range("K90").Value = 2
$ws.Range("N90").Value = 2
$ws.Range("Q90").Value = 2

# Row 91
$ws.Range("E87").Copy()
$ws.Range("H91").PasteSpecial(-4122)
$ws.Range("H91").Formula = "=E91+1"
$ws.Range("K89").Copy()
$ws.Range("K91").PasteSpecial(-4122)
$ws.Range("K91").Value = 2
$ws.Range("N91").Value = 2
$ws.Range("Q91").Value = 2

# ---------------------------------------------------------------------------
# 2) Minor border fix on D73 (now gets the regular thin-bottom border used by
#    the equivalent cell in the following table, D87).
# ---------------------------------------------------------------------------
$ws.Range("D87").Copy()
$ws.Range("D73").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) New underlined helper cell A97 (also bumps sheet dimension from B2:R214
#    to A2:R214 and refreshes the cached row spans for rows 97-112).
# ---------------------------------------------------------------------------
$ws.Range("A97").Value = ""
$ws.Range("A97").Font.Underline = 2

# ---------------------------------------------------------------------------
# 4) Summary rows 111-114: replace the MIN(...) formulas that pulled from the
#    5 weekly tables with direct references to each table's own MIN row.
# ---------------------------------------------------------------------------
$ws.Range("D111").Formula = "=F`$74"
$ws.Range("E111").Formula = "=`$I`$74"
$ws.Range("F111").Formula = "=`$L`$74"
$ws.Range("G111").Formula = "=`$O`$74"
$ws.Range("H111").Formula = "=`$R`$74"

$ws.Range("D112").Formula = "=F83"
$ws.Range("E112").Formula = "=`$I`$83"
$ws.Range("F112").Formula = "=`$L`$83"
$ws.Range("G112").Formula = "=`$O`$83"
$ws.Range("H112").Formula = "=`$R`$83"

$ws.Range("D113").Formula = "=F92"
$ws.Range("E113").Formula = "=`$I`$92"
$ws.Range("F113").Formula = "=`$L`$92"
$ws.Range("G113").Formula = "=`$O`$92"
$ws.Range("H113").Formula = "=`$R`$92"

$ws.Range("E113").Copy()
$ws.Range("E114:H114").PasteSpecial(-4122)
$ws.Range("D114").Formula = "=F101"
$ws.Range("E114").Formula = "=`$I`$101"
$ws.Range("F114").Formula = "=`$L`$101"
$ws.Range("G114").Formula = "=`$O`$101"
$ws.Range("H114").Formula = "=`$R`$101"

# ---------------------------------------------------------------------------
# 5) Sheet view: split the window roughly where the author had it and park
#    selections close to the recorded ones. (Cosmetic only - no data impact.)
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.SplitRow = 16
$win.Split = $true
$ws.Range("G12").Select()
$win.ScrollRow = 84
$ws.Range("E90").Select()

$excel.CutCopyMode = $false
